# Sync attendance_reports: reorder the comma-separated "Recorded By" (column G)
# names on the "Session Analysis Results" sheet so that the special "System"
# entry is sorted to the front of the list (immediately after a leading
# lowercase "system" token, if one is present); lists that contain no
# "System" token are sorted alphabetically instead.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Reorder-RecordedBy {
    param($val)

    if ($val -eq $null) {
        return $val
    }

    $rawParts = $val.Split(",")
    if ($rawParts.Count -lt 2) {
        return $val
    }

    $parts = @()
    foreach ($p in $rawParts) {
        $parts += $p.Trim()
    }

    # locate an exact (case-sensitive) "System" token
    $idx = -1
    for ($i = 0; $i -lt $parts.Count; $i++) {
        if ($parts[$i].Equals("System")) {
            $idx = $i
        }
    }

    if ($idx -ge 0) {
        # target slot: right after a leading lowercase "system", else the front
        $target = 0
        if ($parts[0].Equals("system")) {
            $target = 1
        }

        if ($idx -eq $target) {
            return [string]::Join(", ", $parts)
        }

        $sysVal = $parts[$idx]
        $rest = @()
        for ($i = 0; $i -lt $parts.Count; $i++) {
            if ($i -ne $idx) { $rest += $parts[$i] }
        }

        $final = @()
        for ($i = 0; $i -lt $rest.Count; $i++) {
            if ($i -eq $target) { $final += $sysVal }
            $final += $rest[$i]
        }
        if ($target -ge $rest.Count) { $final += $sysVal }

        return [string]::Join(", ", $final)
    }
    else {
        $sortedParts = $parts | Sort-Object
        return [string]::Join(", ", $sortedParts)
    }
}

$used = $ws.UsedRange
$firstRow = $used.Row
$lastRow = $firstRow + $used.Rows.Count - 1

$changedCount = 0
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $orig = $cell.Value2
    if ($orig -eq $null -or $orig -eq "") {
        continue
    }
    $new = Reorder-RecordedBy $orig
    if ($new -ne $orig) {
        $cell.Value = $new
        $changedCount++
    }
}

Write-Host "Rows updated: $changedCount"
